# Update leve-flip profitability figures (currentAveragePrice / LevePrice / LeveProfit
# columns H-N) across all eight job sheets, refreshed from the latest Universalis
# market snapshot pulled by the scheduled runner.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 2161.25
$ws.Range("I19").Value = 1825
$ws.Range("J19").Value = 2497.5
$ws.Range("K19").Value = 1825
$ws.Range("L19").Value = 2497.5
$ws.Range("M19").Value = -1650
$ws.Range("N19").Value = -2847.5
$ws.Range("H31").Value = 850
$ws.Range("I31").Value = 850
$ws.Range("K31").Value = 2550
$ws.Range("M31").Value = -2320
$ws.Range("H40").Value = 5000.6924
$ws.Range("I40").Value = 3340
$ws.Range("J40").Value = 6038.625
$ws.Range("K40").Value = 3340
$ws.Range("L40").Value = 6038.625
$ws.Range("M40").Value = -3165
$ws.Range("N40").Value = -6388.625
$ws.Range("H51").Value = 6926.846
$ws.Range("J51").Value = 7390.909
$ws.Range("L51").Value = 7390.909
$ws.Range("N51").Value = -8358.909
$ws.Range("H64").Value = 17001
$ws.Range("I64").Value = 0
$ws.Range("J64").Value = 17001
$ws.Range("K64").Value = 0
$ws.Range("L64").Value = 17001
$ws.Range("M64").ClearContents()
$ws.Range("N64").Value = -17497
$ws.Range("H67").Value = 17001
$ws.Range("I67").Value = 0
$ws.Range("J67").Value = 17001
$ws.Range("K67").Value = 0
$ws.Range("L67").Value = 17001
$ws.Range("M67").ClearContents()
$ws.Range("N67").Value = -18717
$ws.Range("H88").Value = 204130.73
$ws.Range("I88").Value = 431132
$ws.Range("J88").Value = 5504.625
$ws.Range("K88").Value = 431132
$ws.Range("L88").Value = 5504.625
$ws.Range("M88").Value = -430726
$ws.Range("N88").Value = -6316.625
$ws.Range("H91").Value = 204130.73
$ws.Range("I91").Value = 431132
$ws.Range("J91").Value = 5504.625
$ws.Range("K91").Value = 431132
$ws.Range("L91").Value = 5504.625
$ws.Range("M91").Value = -429728
$ws.Range("N91").Value = -8312.625
$ws.Range("H96").Value = 691.8946999999999
$ws.Range("I96").Value = 1357
$ws.Range("J96").Value = 384.92307
$ws.Range("K96").Value = 4071
$ws.Range("L96").Value = 1154.76921
$ws.Range("M96").Value = -2698
$ws.Range("N96").Value = -3900.76921
$ws.Range("H98").Value = 201761.3
$ws.Range("I98").Value = 882.55554
$ws.Range("J98").Value = 2009670
$ws.Range("K98").Value = 882.55554
$ws.Range("L98").Value = 2009670
$ws.Range("M98").Value = 615.44446
$ws.Range("N98").Value = -2012666
$ws.Range("H99").Value = 3128.1
$ws.Range("I99").Value = 3398.8
$ws.Range("J99").Value = 2857.4
$ws.Range("K99").Value = 10196.4
$ws.Range("L99").Value = 8572.200000000001
$ws.Range("M99").Value = -8698.400000000001
$ws.Range("N99").Value = -11568.2
$ws.Range("H100").Value = 6936.2856
$ws.Range("I100").Value = 4000
$ws.Range("J100").Value = 7425.6665
$ws.Range("K100").Value = 4000
$ws.Range("L100").Value = 7425.6665
$ws.Range("M100").Value = -3459
$ws.Range("N100").Value = -8507.666499999999
$ws.Range("H101").Value = 2400.4167
$ws.Range("I101").Value = 1256.1111
$ws.Range("K101").Value = 3768.3333
$ws.Range("M101").Value = -2146.3333
$ws.Range("H112").Value = 2108.4644
$ws.Range("J112").Value = 2108.4644
$ws.Range("L112").Value = 6325.3932
$ws.Range("N112").Value = -8541.393199999999
$ws.Range("H121").Value = 1147
$ws.Range("J121").Value = 1147
$ws.Range("L121").Value = 3441
$ws.Range("N121").Value = -6935
$ws.Range("H122").Value = 201761.3
$ws.Range("I122").Value = 882.55554
$ws.Range("J122").Value = 2009670
$ws.Range("K122").Value = 2647.66662
$ws.Range("L122").Value = 6029010
$ws.Range("M122").Value = -197.66662
$ws.Range("N122").Value = -6033910
$ws.Range("H123").Value = 66984
$ws.Range("J123").Value = 66984
$ws.Range("L123").Value = 66984
$ws.Range("N123").Value = -76784
$ws.Range("H127").Value = 2184.5
$ws.Range("I127").Value = 1497.6666
$ws.Range("J127").Value = 4245
$ws.Range("K127").Value = 4492.9998
$ws.Range("L127").Value = 12735
$ws.Range("M127").Value = 467.0002000000004
$ws.Range("N127").Value = -22655
$ws.Range("H131").Value = 5471.7407
$ws.Range("I131").Value = 3058.875
$ws.Range("K131").Value = 9176.625
$ws.Range("M131").Value = -4136.625
$ws.Range("H132").Value = 1828.75
$ws.Range("I132").Value = 1732.1632
$ws.Range("J132").Value = 3406.3333
$ws.Range("K132").Value = 5196.4896
$ws.Range("L132").Value = 10218.9999
$ws.Range("M132").Value = -2666.4896
$ws.Range("N132").Value = -15278.9999
$ws.Range("H137").Value = 13892091
$ws.Range("I137").Value = 52634424
$ws.Range("K137").Value = 157903272
$ws.Range("M137").Value = -157900722
$ws.Range("H138").Value = 4450.18
$ws.Range("I138").Value = 1953.3334
$ws.Range("K138").Value = 5860.0002
$ws.Range("M138").Value = -720.0002000000004
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 6492.75
$ws.Range("I2").Value = 1697.3334
$ws.Range("J2").Value = 13685.875
$ws.Range("K2").Value = 1697.3334
$ws.Range("L2").Value = 13685.875
$ws.Range("M2").Value = -1584.3334
$ws.Range("N2").Value = -13911.875
$ws.Range("H5").Value = 283.9375
$ws.Range("I5").Value = 232
$ws.Range("J5").Value = 370.5
$ws.Range("K5").Value = 232
$ws.Range("L5").Value = 370.5
$ws.Range("M5").Value = -120
$ws.Range("N5").Value = -594.5
$ws.Range("H32").Value = 9439.261
$ws.Range("I32").Value = 8186.1816
$ws.Range("K32").Value = 8186.1816
$ws.Range("M32").Value = -7899.1816
$ws.Range("H45").Value = 3593.5293
$ws.Range("I45").Value = 3208.7307
$ws.Range("J45").Value = 4844.125
$ws.Range("K45").Value = 3208.7307
$ws.Range("L45").Value = 4844.125
$ws.Range("M45").Value = -2831.7307
$ws.Range("N45").Value = -5598.125
$ws.Range("H63").Value = 10000
$ws.Range("I63").Value = 0
$ws.Range("K63").Value = 0
$ws.Range("M63").ClearContents()
$ws.Range("H66").Value = 10000
$ws.Range("I66").Value = 0
$ws.Range("K66").Value = 0
$ws.Range("M66").ClearContents()
$ws.Range("H74").Value = 6340.5713
$ws.Range("I74").Value = 4374
$ws.Range("J74").Value = 7815.5
$ws.Range("K74").Value = 4374
$ws.Range("L74").Value = 7815.5
$ws.Range("M74").Value = -3500
$ws.Range("N74").Value = -9563.5
$ws.Range("H77").Value = 6340.5713
$ws.Range("I77").Value = 4374
$ws.Range("J77").Value = 7815.5
$ws.Range("K77").Value = 21870
$ws.Range("L77").Value = 39077.5
$ws.Range("M77").Value = -17502
$ws.Range("N77").Value = -47813.5
$ws.Range("H88").Value = 20818.4
$ws.Range("J88").Value = 19991.143
$ws.Range("L88").Value = 19991.143
$ws.Range("N88").Value = -20803.143
$ws.Range("H91").Value = 20818.4
$ws.Range("J91").Value = 19991.143
$ws.Range("L91").Value = 19991.143
$ws.Range("N91").Value = -22799.143
$ws.Range("H97").Value = 924.7778
$ws.Range("I97").Value = 811.9167
$ws.Range("J97").Value = 1150.5
$ws.Range("K97").Value = 811.9167
$ws.Range("L97").Value = 1150.5
$ws.Range("M97").Value = -315.9167
$ws.Range("N97").Value = -2142.5
$ws.Range("H110").Value = 4119
$ws.Range("I110").Value = 3052.2222
$ws.Range("J110").Value = 6862.143
$ws.Range("K110").Value = 3052.2222
$ws.Range("L110").Value = 6862.143
$ws.Range("M110").Value = -1007.2222
$ws.Range("N110").Value = -10952.143
$ws.Range("H116").Value = 6492.75
$ws.Range("I116").Value = 1697.3334
$ws.Range("J116").Value = 13685.875
$ws.Range("K116").Value = 1697.3334
$ws.Range("L116").Value = 13685.875
$ws.Range("M116").Value = 596.6666
$ws.Range("N116").Value = -18273.875
$ws.Range("H121").Value = 41662.332
$ws.Range("J121").Value = 41662.332
$ws.Range("L121").Value = 41662.332
$ws.Range("N121").Value = -45156.332
$ws.Range("H122").Value = 4824.0625
$ws.Range("I122").Value = 4962.125
$ws.Range("J122").Value = 4686
$ws.Range("K122").Value = 14886.375
$ws.Range("L122").Value = 14058
$ws.Range("M122").Value = -12436.375
$ws.Range("N122").Value = -18958
$ws.Range("H132").Value = 4128.4062
$ws.Range("I132").Value = 3610.5715
$ws.Range("K132").Value = 10831.7145
$ws.Range("M132").Value = -8301.7145
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 6492.75
$ws.Range("I3").Value = 1697.3334
$ws.Range("J3").Value = 13685.875
$ws.Range("K3").Value = 1697.3334
$ws.Range("L3").Value = 13685.875
$ws.Range("M3").Value = -1583.3334
$ws.Range("N3").Value = -13913.875
$ws.Range("H4").Value = 283.9375
$ws.Range("I4").Value = 232
$ws.Range("J4").Value = 370.5
$ws.Range("K4").Value = 232
$ws.Range("L4").Value = 370.5
$ws.Range("M4").Value = -117
$ws.Range("N4").Value = -600.5
$ws.Range("H26").Value = 43416.5
$ws.Range("I26").Value = 28875
$ws.Range("K26").Value = 28875
$ws.Range("M26").Value = -28583
$ws.Range("H64").Value = 1259.875
$ws.Range("J64").Value = 1297
$ws.Range("L64").Value = 1297
$ws.Range("N64").Value = -1747
$ws.Range("H67").Value = 1259.875
$ws.Range("J67").Value = 1297
$ws.Range("L67").Value = 1297
$ws.Range("N67").Value = -2857
$ws.Range("H86").Value = 4334.2173
$ws.Range("I86").Value = 3452.0667
$ws.Range("J86").Value = 5988.25
$ws.Range("K86").Value = 3452.0667
$ws.Range("L86").Value = 5988.25
$ws.Range("M86").Value = -2329.0667
$ws.Range("N86").Value = -8234.25
$ws.Range("H89").Value = 4334.2173
$ws.Range("I89").Value = 3452.0667
$ws.Range("J89").Value = 5988.25
$ws.Range("K89").Value = 17260.3335
$ws.Range("L89").Value = 29941.25
$ws.Range("M89").Value = -11644.3335
$ws.Range("N89").Value = -41173.25
$ws.Range("H94").Value = 1774
$ws.Range("I94").Value = 1769.5
$ws.Range("J94").Value = 1798
$ws.Range("K94").Value = 1769.5
$ws.Range("L94").Value = 1798
$ws.Range("M94").Value = -1318.5
$ws.Range("N94").Value = -2700
$ws.Range("H133").Value = 74966.664
$ws.Range("J133").Value = 74966.664
$ws.Range("L133").Value = 74966.664
$ws.Range("N133").Value = -85086.664
$ws.Range("H134").Value = 4304.3477
$ws.Range("I134").Value = 3750
$ws.Range("K134").Value = 11250
$ws.Range("M134").Value = -8715
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H21").Value = 2550
$ws.Range("J21").Value = 2550
$ws.Range("L21").Value = 2550
$ws.Range("N21").Value = -3020
$ws.Range("H22").Value = 7044.222
$ws.Range("I22").Value = 2233
$ws.Range("J22").Value = 16666.666
$ws.Range("K22").Value = 2233
$ws.Range("L22").Value = 16666.666
$ws.Range("M22").Value = -1883
$ws.Range("N22").Value = -17366.666
$ws.Range("H31").Value = 37190.324
$ws.Range("I31").Value = 8059.2383
$ws.Range("J31").Value = 84248.234
$ws.Range("K31").Value = 8059.2383
$ws.Range("L31").Value = 84248.234
$ws.Range("M31").Value = -7764.2383
$ws.Range("N31").Value = -84838.234
$ws.Range("H34").Value = 37190.324
$ws.Range("I34").Value = 8059.2383
$ws.Range("J34").Value = 84248.234
$ws.Range("K34").Value = 8059.2383
$ws.Range("L34").Value = 84248.234
$ws.Range("M34").Value = -7857.2383
$ws.Range("N34").Value = -84652.234
$ws.Range("H58").Value = 7613
$ws.Range("I58").Value = 2012
$ws.Range("J58").Value = 8631.362999999999
$ws.Range("K58").Value = 2012
$ws.Range("L58").Value = 8631.362999999999
$ws.Range("M58").Value = -1809
$ws.Range("N58").Value = -9037.362999999999
$ws.Range("H82").Value = 34981
$ws.Range("J82").Value = 34981
$ws.Range("L82").Value = 34981
$ws.Range("N82").Value = -35703
$ws.Range("H85").Value = 34981
$ws.Range("J85").Value = 34981
$ws.Range("L85").Value = 34981
$ws.Range("N85").Value = -37477
$ws.Range("H97").Value = 59999.332
$ws.Range("I97").Value = 49999
$ws.Range("J97").Value = 64999.5
$ws.Range("K97").Value = 49999
$ws.Range("L97").Value = 64999.5
$ws.Range("M97").Value = -49008
$ws.Range("N97").Value = -66981.5
$ws.Range("H99").Value = 4599.2
$ws.Range("I99").Value = 3999
$ws.Range("K99").Value = 3999
$ws.Range("M99").Value = -2501
$ws.Range("H122").Value = 4347.08
$ws.Range("I122").Value = 1510.8667
$ws.Range("J122").Value = 8601.4
$ws.Range("K122").Value = 4532.6001
$ws.Range("L122").Value = 25804.2
$ws.Range("M122").Value = -2082.6001
$ws.Range("N122").Value = -30704.2
$ws.Range("H126").Value = 4599.2
$ws.Range("I126").Value = 3999
$ws.Range("K126").Value = 11997
$ws.Range("M126").Value = -9527
$ws.Range("H132").Value = 3107.6394
$ws.Range("I132").Value = 2707.6843
$ws.Range("K132").Value = 8123.0529
$ws.Range("M132").Value = -5593.0529
$ws.Range("H134").Value = 3588.4443
$ws.Range("I134").Value = 2238.5334
$ws.Range("J134").Value = 10338
$ws.Range("K134").Value = 6715.600199999999
$ws.Range("L134").Value = 31014
$ws.Range("M134").Value = -4180.600199999999
$ws.Range("N134").Value = -36084
$ws.Range("H136").Value = 7613
$ws.Range("I136").Value = 2012
$ws.Range("J136").Value = 8631.362999999999
$ws.Range("K136").Value = 6036
$ws.Range("L136").Value = 25894.089
$ws.Range("M136").Value = -3486
$ws.Range("N136").Value = -30994.089
$ws.Range("H141").Value = 345786
$ws.Range("J141").Value = 361665.78
$ws.Range("L141").Value = 361665.78
$ws.Range("N141").Value = -372025.78
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 123159.08
$ws.Range("I2").Value = 64.833336
$ws.Range("J2").Value = 228668.42
$ws.Range("K2").Value = 389.000016
$ws.Range("L2").Value = 1372010.52
$ws.Range("M2").Value = -276.000016
$ws.Range("N2").Value = -1372236.52
$ws.Range("H14").Value = 118
$ws.Range("I14").Value = 118
$ws.Range("K14").Value = 354
$ws.Range("M14").Value = -181
$ws.Range("H34").Value = 4172104.8
$ws.Range("I34").Value = 5925.6
$ws.Range("J34").Value = 25003000
$ws.Range("K34").Value = 17776.8
$ws.Range("L34").Value = 75009000
$ws.Range("M34").Value = -17692.8
$ws.Range("N34").Value = -75009168
$ws.Range("H39").Value = 2530
$ws.Range("J39").Value = 2466.6667
$ws.Range("L39").Value = 7400.000100000001
$ws.Range("N39").Value = -7988.000100000001
$ws.Range("H80").Value = 7249.25
$ws.Range("J80").Value = 6799.6
$ws.Range("L80").Value = 20398.8
$ws.Range("N80").Value = -22270.8
$ws.Range("H83").Value = 7249.25
$ws.Range("J83").Value = 6799.6
$ws.Range("L83").Value = 61196.4
$ws.Range("N83").Value = -70556.39999999999
$ws.Range("H92").Value = 922.1177
$ws.Range("I92").Value = 572.5
$ws.Range("J92").Value = 1029.6923
$ws.Range("K92").Value = 1717.5
$ws.Range("L92").Value = 3089.0769
$ws.Range("M92").Value = -469.5
$ws.Range("N92").Value = -5585.0769
$ws.Range("H122").Value = 2508.3872
$ws.Range("I122").Value = 823.6667
$ws.Range("J122").Value = 3572.4211
$ws.Range("K122").Value = 7413.0003
$ws.Range("L122").Value = 32151.7899
$ws.Range("M122").Value = -4963.0003
$ws.Range("N122").Value = -37051.7899
$ws.Range("H131").Value = 18792506
$ws.Range("I131").Value = 31251542
$ws.Range("J131").Value = 14958957
$ws.Range("K131").Value = 93754626
$ws.Range("L131").Value = 44876871
$ws.Range("M131").Value = -93749586
$ws.Range("N131").Value = -44886951
$ws.Range("H132").Value = 4976.7393
$ws.Range("I132").Value = 4349
$ws.Range("K132").Value = 39141
$ws.Range("M132").Value = -36611
$ws.Range("H141").Value = 3563.0435
$ws.Range("I141").Value = 1273.1666
$ws.Range("K141").Value = 3819.4998
$ws.Range("M141").Value = 1360.5002
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H11").Value = 10375000
$ws.Range("I11").Value = 10000000
$ws.Range("K11").Value = 10000000
$ws.Range("M11").Value = -9999861
$ws.Range("H34").Value = 45000
$ws.Range("J34").Value = 45000
$ws.Range("L34").Value = 45000
$ws.Range("N34").Value = -45536
$ws.Range("H62").Value = 44999
$ws.Range("I62").Value = 44999
$ws.Range("K62").Value = 44999
$ws.Range("M62").Value = -44313
$ws.Range("H65").Value = 44999
$ws.Range("I65").Value = 44999
$ws.Range("K65").Value = 134997
$ws.Range("M65").Value = -131565
$ws.Range("H70").Value = 10463.789
$ws.Range("I70").Value = 9620.429
$ws.Range("J70").Value = 12825.2
$ws.Range("K70").Value = 9620.429
$ws.Range("L70").Value = 12825.2
$ws.Range("M70").Value = -9350.429
$ws.Range("N70").Value = -13365.2
$ws.Range("H73").Value = 10463.789
$ws.Range("I73").Value = 9620.429
$ws.Range("J73").Value = 12825.2
$ws.Range("K73").Value = 9620.429
$ws.Range("L73").Value = 12825.2
$ws.Range("M73").Value = -8684.429
$ws.Range("N73").Value = -14697.2
$ws.Range("H76").Value = 45000
$ws.Range("J76").Value = 45000
$ws.Range("L76").Value = 45000
$ws.Range("N76").Value = -45630
$ws.Range("H79").Value = 45000
$ws.Range("J79").Value = 45000
$ws.Range("L79").Value = 45000
$ws.Range("N79").Value = -47184
$ws.Range("H97").Value = 1676.35
$ws.Range("I97").Value = 732.96
$ws.Range("J97").Value = 3248.6667
$ws.Range("K97").Value = 732.96
$ws.Range("L97").Value = 3248.6667
$ws.Range("M97").Value = -236.96
$ws.Range("N97").Value = -4240.6667
$ws.Range("H102").Value = 1576.1708
$ws.Range("I102").Value = 512.8125
$ws.Range("K102").Value = 512.8125
$ws.Range("M102").Value = 1109.1875
$ws.Range("H113").Value = 3921.7222
$ws.Range("I113").Value = 3137
$ws.Range("J113").Value = 4421.091
$ws.Range("K113").Value = 3137
$ws.Range("L113").Value = 4421.091
$ws.Range("M113").Value = -967
$ws.Range("N113").Value = -8761.091
$ws.Range("H122").Value = 6328.55
$ws.Range("I122").Value = 4714.7856
$ws.Range("J122").Value = 10094
$ws.Range("K122").Value = 14144.3568
$ws.Range("L122").Value = 30282
$ws.Range("M122").Value = -11694.3568
$ws.Range("N122").Value = -35182
$ws.Range("H132").Value = 3494.75
$ws.Range("I132").Value = 2092.3
$ws.Range("K132").Value = 6276.900000000001
$ws.Range("M132").Value = -3746.900000000001
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 8940.817999999999
$ws.Range("I7").Value = 6394
$ws.Range("J7").Value = 17600
$ws.Range("K7").Value = 6394
$ws.Range("L7").Value = 17600
$ws.Range("M7").Value = -6282
$ws.Range("N7").Value = -17824
$ws.Range("H22").Value = 11937.125
$ws.Range("I22").Value = 1500
$ws.Range("J22").Value = 13428.143
$ws.Range("K22").Value = 1500
$ws.Range("L22").Value = 13428.143
$ws.Range("M22").Value = -1205
$ws.Range("N22").Value = -14018.143
$ws.Range("H27").Value = 11937.125
$ws.Range("I27").Value = 1500
$ws.Range("J27").Value = 13428.143
$ws.Range("K27").Value = 1500
$ws.Range("L27").Value = 13428.143
$ws.Range("M27").Value = -1393
$ws.Range("N27").Value = -13642.143
$ws.Range("H40").Value = 9187.423000000001
$ws.Range("I40").Value = 10948.182
$ws.Range("J40").Value = 7896.2
$ws.Range("K40").Value = 10948.182
$ws.Range("L40").Value = 7896.2
$ws.Range("M40").Value = -10812.182
$ws.Range("N40").Value = -8168.2
$ws.Range("H46").Value = 4650.3
$ws.Range("J46").Value = 4944.6665
$ws.Range("L46").Value = 4944.6665
$ws.Range("N46").Value = -5320.6665
$ws.Range("H68").Value = 6053.2144
$ws.Range("I68").Value = 3665.75
$ws.Range("K68").Value = 3665.75
$ws.Range("M68").Value = -2916.75
$ws.Range("H71").Value = 6053.2144
$ws.Range("I71").Value = 3665.75
$ws.Range("K71").Value = 18328.75
$ws.Range("M71").Value = -14584.75
$ws.Range("H82").Value = 5152.9565
$ws.Range("I82").Value = 3424.2
$ws.Range("J82").Value = 6482.769
$ws.Range("K82").Value = 3424.2
$ws.Range("L82").Value = 6482.769
$ws.Range("M82").Value = -3063.2
$ws.Range("N82").Value = -7204.769
$ws.Range("H85").Value = 5152.9565
$ws.Range("I85").Value = 3424.2
$ws.Range("J85").Value = 6482.769
$ws.Range("K85").Value = 3424.2
$ws.Range("L85").Value = 6482.769
$ws.Range("M85").Value = -2176.2
$ws.Range("N85").Value = -8978.769
$ws.Range("H87").Value = 39999
$ws.Range("I87").Value = 39999
$ws.Range("K87").Value = 39999
$ws.Range("M87").Value = -38876
$ws.Range("H90").Value = 39999
$ws.Range("I90").Value = 39999
$ws.Range("K90").Value = 119997
$ws.Range("M90").Value = -114381
$ws.Range("H93").Value = 2377.3333
$ws.Range("I93").Value = 2756.8
$ws.Range("J93").Value = 480
$ws.Range("K93").Value = 2756.8
$ws.Range("L93").Value = 480
$ws.Range("M93").Value = -1508.8
$ws.Range("N93").Value = -2976
$ws.Range("H122").Value = 257522.69
$ws.Range("I122").Value = 450820.88
$ws.Range("K122").Value = 1352462.64
$ws.Range("M122").Value = -1350012.64
$ws.Range("H126").Value = 8940.817999999999
$ws.Range("I126").Value = 6394
$ws.Range("J126").Value = 17600
$ws.Range("K126").Value = 19182
$ws.Range("L126").Value = 52800
$ws.Range("M126").Value = -16712
$ws.Range("N126").Value = -57740
$ws.Range("H132").Value = 3933.4773
$ws.Range("I132").Value = 3528.1052
$ws.Range("J132").Value = 6500.8335
$ws.Range("K132").Value = 10584.3156
$ws.Range("L132").Value = 19502.5005
$ws.Range("M132").Value = -8054.3156
$ws.Range("N132").Value = -24562.5005
$ws.Range("H136").Value = 4552.1055
$ws.Range("I136").Value = 3198.5
$ws.Range("J136").Value = 6056.1113
$ws.Range("K136").Value = 9595.5
$ws.Range("L136").Value = 18168.3339
$ws.Range("M136").Value = -7045.5
$ws.Range("N136").Value = -23268.3339
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H15").Value = 39999.5
$ws.Range("J15").Value = 39999.5
$ws.Range("L15").Value = 39999.5
$ws.Range("N15").Value = -40575.5
$ws.Range("H50").Value = 0
$ws.Range("J50").Value = 0
$ws.Range("L50").Value = 0
$ws.Range("N50").ClearContents()
$ws.Range("H62").Value = 7088.4614
$ws.Range("I62").Value = 4216.3335
$ws.Range("J62").Value = 7950.1
$ws.Range("K62").Value = 4216.3335
$ws.Range("L62").Value = 7950.1
$ws.Range("M62").Value = -3592.3335
$ws.Range("N62").Value = -9198.1
$ws.Range("H65").Value = 7088.4614
$ws.Range("I65").Value = 4216.3335
$ws.Range("J65").Value = 7950.1
$ws.Range("K65").Value = 21081.6675
$ws.Range("L65").Value = 39750.5
$ws.Range("M65").Value = -17961.6675
$ws.Range("N65").Value = -45990.5
$ws.Range("H75").Value = 49331.332
$ws.Range("I75").Value = 49331.332
$ws.Range("K75").Value = 49331.332
$ws.Range("M75").Value = -48395.332
$ws.Range("H78").Value = 49331.332
$ws.Range("I78").Value = 49331.332
$ws.Range("K78").Value = 147993.996
$ws.Range("M78").Value = -143313.996
$ws.Range("H81").Value = 7704.8667
$ws.Range("I81").Value = 6853.125
$ws.Range("K81").Value = 13706.25
$ws.Range("M81").Value = -12645.25
$ws.Range("H84").Value = 7704.8667
$ws.Range("I84").Value = 6853.125
$ws.Range("K84").Value = 68531.25
$ws.Range("M84").Value = -63227.25
$ws.Range("H86").Value = 47625
$ws.Range("I86").Value = 0
$ws.Range("K86").Value = 0
$ws.Range("M86").ClearContents()
$ws.Range("H89").Value = 47625
$ws.Range("I89").Value = 0
$ws.Range("K89").Value = 0
$ws.Range("M89").ClearContents()
$ws.Range("H96").Value = 2860.4614
$ws.Range("I96").Value = 1839.6
$ws.Range("J96").Value = 3498.5
$ws.Range("K96").Value = 1839.6
$ws.Range("L96").Value = 3498.5
$ws.Range("M96").Value = -466.5999999999999
$ws.Range("N96").Value = -6244.5
$ws.Range("H100").Value = 973.9355
$ws.Range("I100").Value = 781.0909
$ws.Range("J100").Value = 1445.3334
$ws.Range("K100").Value = 1562.1818
$ws.Range("L100").Value = 2890.6668
$ws.Range("M100").Value = -1021.1818
$ws.Range("N100").Value = -3972.6668
$ws.Range("H113").Value = 282.5
$ws.Range("I113").Value = 305.70834
$ws.Range("J113").Value = 143.25
$ws.Range("K113").Value = 917.1250200000001
$ws.Range("L113").Value = 429.75
$ws.Range("M113").Value = 1252.87498
$ws.Range("N113").Value = -4769.75
$ws.Range("H122").Value = 2545.1538
$ws.Range("I122").Value = 1392.8889
$ws.Range("J122").Value = 5137.75
$ws.Range("K122").Value = 4178.6667
$ws.Range("L122").Value = 15413.25
$ws.Range("M122").Value = -1728.6667
$ws.Range("N122").Value = -20313.25
$ws.Range("H126").Value = 3079.0625
$ws.Range("I126").Value = 1723.3334
$ws.Range("K126").Value = 5170.0002
$ws.Range("M126").Value = -2700.0002
$ws.Range("H132").Value = 6187.0713
$ws.Range("I132").Value = 3735
$ws.Range("J132").Value = 10600.8
$ws.Range("K132").Value = 11205
$ws.Range("L132").Value = 31802.4
$ws.Range("M132").Value = -8675
$ws.Range("N132").Value = -36862.39999999999
$ws.Range("H136").Value = 4543.857
$ws.Range("I136").Value = 2070.7273
$ws.Range("K136").Value = 6212.1819
$ws.Range("M136").Value = -3662.1819
